$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.0681
$ws.Range("D6").Value = -8.4596
$ws.Range("A14").Value = -21.81880000000001
$ws.Range("D18").Value = -8.5487
$ws.Range("D19").Value = -8.429799999999998
$ws.Range("A21").Value = -20.04919999999997
$ws.Range("B22").Value = 9.487299999999992
$ws.Range("A23").Value = -20.43219999999998
$ws.Range("B24").Value = 5.578500000000002
$ws.Range("A25").Value = -21.65069999999998
$ws.Range("A26").Value = -21.00729999999997
$ws.Range("B28").Value = 5.041900000000002
$ws.Range("A29").Value = -21.26589999999998
$ws.Range("B36").Value = 9.450300000000009
$ws.Range("D44").Value = -6.888499999999994
$ws.Range("B45").Value = 4.7818
$ws.Range("D47").Value = -7.2823
$ws.Range("B48").Value = 7.981500000000004
$ws.Range("B49").Value = 5.752499999999998
$ws.Range("D51").Value = -8.250599999999997
$ws.Range("B52").Value = 6.130100000000001
$ws.Range("A53").Value = -21.75739999999998
$ws.Range("B53").Value = 5.800599999999999
$ws.Range("B54").Value = 4.815200000000001
$ws.Range("D55").Value = -9.0063
$ws.Range("A57").Value = -22.45470000000002
$ws.Range("D57").Value = -8.168899999999997
$ws.Range("A59").Value = -22.26449999999999
$ws.Range("D64").Value = -7.660199999999991
$ws.Range("A69").Value = -21.63099999999999
$ws.Range("B70").Value = 7.652200000000004
$ws.Range("A79").Value = -20.43950000000001
$ws.Range("D80").Value = -7.620000000000004
$ws.Range("A83").Value = -21.77339999999999
$ws.Range("B86").Value = 4.8565
$ws.Range("B87").Value = 5.420700000000003
$ws.Range("B89").Value = 4.477599999999999
$ws.Range("A91").Value = -20.65059999999999
$ws.Range("D92").Value = -6.698300000000004
$ws.Range("A93").Value = -21.41760000000002
$ws.Range("D94").Value = -6.589900000000003
$ws.Range("D96").Value = -8.484100000000005
$ws.Range("B101").Value = 5.1296
$ws.Range("D101").Value = -8.137099999999997
$ws.Range("A103").Value = -21.75799999999999
